# "Forms SOPs 1 & 2"
# Rename the Service/Request form sheet from F-SW-SD-01 to F-SW-FR-01,
# change the form title text accordingly, update the Print_Area defined
# name to follow the new sheet name, and refresh the selection/zoom on
# the (now renamed) sheet to match the saved view state.

$wb = $excel.ActiveWorkbook

# --- Rename the main sheet ---------------------------------------------
$ws = $wb.Worksheets.Item("F-SW-SD-01")
$ws.Name = "F-SW-FR-01"

# --- Update the big heading in A3 from "SERVICE REQUEST FORM" to
#     "REQUEST FORM" -------------------------------------------------
$ws.Range("A3").Value = "REQUEST FORM"

# --- Keep the Print_Area defined name pointing at the renamed sheet ----
foreach ($n in $wb.Names) {
    if ($n.Name -eq "F-SW-FR-01!Print_Area") {
        $n.RefersTo = "='F-SW-FR-01'!`$A`$1:`$F`$22"
    }
}

# --- Update the saved view: 100% zoom and selection on the department
#     drop-down cell (C10:F10) instead of the old title selection -------
$win = $wb.Windows.Item(1)
$win.Zoom = 100
[void]$ws.Range("C10:F10").Select()
